# Refresh the crypto price/volume snapshot (columns D "Price" and
# E "Volume(1h)") for each coin row on the active sheet, matching the
# latest coinranking.com pull from the GitHub Actions job.
$ws = $excel.ActiveWorkbook.ActiveSheet

# row -> @{ D = "<new price text>"; E = "<new volume text>" } (a column is
# omitted when that cell is unchanged for the row).
$updates = [ordered]@{
    2 = @{ D='23.380.27'; E='  +1.57%  ' }
    3 = @{ D='1.627.32'; E='  +2.66%  ' }
    4 = @{ D='0.9974'; E='  -0.67%  ' }
    5 = @{ D='306.98'; E='  +1.95%  ' }
    6 = @{ D='0.9965'; E='  -0.66%  ' }
    7 = @{ D='0.3782'; E='  +0.74%  ' }
    8 = @{ D='53.25'; E='  +5.01%  ' }
    9 = @{ D='0.3662'; E='  +2.36%  ' }
    10 = @{ D='1.278'; E='  +4.93%  ' }
    11 = @{ D='0.08189'; E='  +2.23%  ' }
    12 = @{ D='0.9953'; E='  -0.90%  ' }
    13 = @{ D='23.20'; E='  +6.20%  ' }
    14 = @{ E='  +3.14%  ' }
    15 = @{ D='7.448'; E='  +2.61%  ' }
    16 = @{ E='  +2.93%  ' }
    17 = @{ D='1.624.61'; E='  +2.54%  ' }
    18 = @{ D='94.71'; E='  +2.69%  ' }
    19 = @{ D='0.06922'; E='  +1.76%  ' }
    20 = @{ D='18.36'; E='  +2.83%  ' }
    21 = @{ D='6.579'; E='  +2.28%  ' }
    22 = @{ E='  -0.31%  ' }
    23 = @{ E='  +1.32%  ' }
    24 = @{ D='23.404.79'; E='  +1.67%  ' }
    25 = @{ D='3.142'; E='  +13.03%  ' }
    26 = @{ D='2.414'; E='  +1.79%  ' }
    27 = @{ D='21.39'; E='  +3.09%  ' }
    28 = @{ D='150.62'; E='  +2.08%  ' }
    29 = @{ D='5.273'; E='  +1.15%  ' }
    30 = @{ D='136.38'; E='  +2.82%  ' }
    31 = @{ D='2.421'; E='  +2.53%  ' }
    32 = @{ D='6.893'; E='  +6.06%  ' }
    33 = @{ D='1.799.32'; E='  +1.83%  ' }
    34 = @{ D='0.9739'; E='  +3.77%  ' }
    35 = @{ D='0.02789'; E='  +4.39%  ' }
    36 = @{ E='  +3.89%  ' }
    37 = @{ D='0.07449'; E='  +1.70%  ' }
    38 = @{ D='6.228'; E='  +3.33%  ' }
    39 = @{ E='  +2.54%  ' }
    40 = @{ E='  +0.96%  ' }
    41 = @{ D='1.405'; E='  +4.67%  ' }
    42 = @{ D='0.7158'; E='  +4.35%  ' }
    43 = @{ E='  +7.12%  ' }
    44 = @{ D='16.19'; E='  +8.76%  ' }
    45 = @{ D='0.6615'; E='  +3.48%  ' }
    46 = @{ D='2.359'; E='  +5.38%  ' }
    47 = @{ D='4.034'; E='  +0.97%  ' }
    48 = @{ D='0.9960'; E='  -0.47%  ' }
    49 = @{ D='0.08012'; E='  +1.66%  ' }
    50 = @{ E='  +0.26%  ' }
    51 = @{ D='1.215'; E='  +2.31%  ' }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $colIndex = @{ D = 4; E = 5 }[$col]
        $text = $updates[$row][$col]
        # Price strings like "0.9974" or "23.20" are indistinguishable from
        # plain numbers, so a bare assignment would make Excel parse them as
        # numeric values and silently drop significant trailing/grouping
        # digits. A leading apostrophe forces literal-text entry, exactly as
        # typing '0.9974 into the cell would - the cell keeps storing plain
        # text (matching the source inline-string cells), just flagged
        # quote-prefixed instead of General-numeric.
        if ($text -match '^[+-]?\d+(\.\d+)?$') {
            $ws.Cells.Item($row, $colIndex).Value = "'" + $text
        } else {
            $ws.Cells.Item($row, $colIndex).Value = $text
        }
    }
}
